# Update the "Spatiotemporal" sheet's row-6 simulated data values to reflect
# the updated spatiotemporal model estimates. Row 7 (SUM/AVERAGE) recalculates
# automatically from these.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Spatiotemporal")

$ws.Range("C6").Value = 4039.4268033950002
$ws.Range("D6").Value = 4.6098403848728902
$ws.Range("F6").Value = 0.288468766059159
$ws.Range("H6").Value = 0.53383034567051602
$ws.Range("I6").Value = 0.27702017411214502
$ws.Range("J6").Value = 1.7845357201518699
$ws.Range("K6").Value = 0.66460491603462701
$ws.Range("L6").Value = 0.18605278237171399
$ws.Range("M6").Value = 0.74313567631079402
$ws.Range("N6").Value = 0.49834392802584299
$ws.Range("O6").Value = 0.448185970392319
$ws.Range("P6").Value = 0.89270664352491302
$ws.Range("Q6").Value = 0.66353831843727196
$ws.Range("R6").Value = 0.83172412895627701
$ws.Range("S6").Value = 0.26400088840503799

# Update the current selection/view on the Spatiotemporal sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D16").Select()

# Update the current selection on the Semantic sheet (view-only change).
$ws2 = $wb.Worksheets.Item("Semantic")
$ws2.Activate()
$ws2.Range("M6").Select()

# Re-activate the Spatiotemporal sheet (it was the active tab originally).
$ws.Activate()
